$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("exiobase").Name = "Exiobase"
$wb.Worksheets.Item("german").Name = "Deutsch"
$wb.Worksheets.Item("english").Name = "Englisch"

# Move the "selected tab" marker from raw_material to Englisch (4th sheet),
# and set its active selection to E32.
$wsRaw = $wb.Worksheets.Item("raw_material")
$wsRaw.Select()

$wsEng = $wb.Worksheets.Item("Englisch")
$wsEng.Select()
$wsEng.Range("E32").Select()
